$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.248.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.982.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.96%  "

# Row 6
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("E8").Value = "  +4.91%  "

# Row 9
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.983.96"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.62%  "

# Row 11
$ws.Range("E11").Value = "  +3.48%  "

# Row 12
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.15%  "

# Row 13
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.35%  "

# Row 14
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.79%  "

# Row 15
$ws.Range("E15").Value = "  +3.18%  "

# Row 16
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.475.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.28%  "

# Row 17
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.09%  "

# Row 18
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.988.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.59%  "

# Row 19
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "58.203.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "423.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.11%  "

# Row 21
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.42%  "

# Row 22
$ws.Range("E22").Value = "  +9.19%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.44%  "

# Row 24
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.78%  "

# Row 25
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.68%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.20%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.84%  "

# Row 30
$ws.Range("E30").Value = "  +6.52%  "

# Row 31
$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.38%  "

# Row 32
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0975"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.87%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.58%  "

# Row 35
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.964"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.76%  "

# Row 36
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0737"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +19.87%  "

# Row 37
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.30%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "

# Row 39
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.92%  "

# Row 40
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.14%  "

# Row 41
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "391.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.09%  "

# Row 42
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0351"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.71%  "

# Row 43
$ws.Range("E43").Value = "  +3.13%  "

# Row 44
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.726.65"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.244"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.33%  "

# Row 46
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.28%  "

# Row 48
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.70%  "

# Row 49
$ws.Range("E49").Value = "  +2.89%  "

# Row 50
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.82%  "

# Row 51
$ws.Range("E51").Value = "  +4.46%  "
